# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) values on the zh-cn and de-de
# report sheets to reflect the latest handback report generation times.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-23 19:22:18"
$zhcn.Range("H2").Value = "2016-03-23 19:22:46"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-23 19:22:23"
$dede.Range("H2").Value = "2016-03-23 19:22:53"
